# Update transcript-related questions to be more open-ended, removing
# specific topic assumptions (credit quality, acquisitions, etc.)

$wb = $excel.ActiveWorkbook

# --- "Landing Page" sheet ---
$ws = $wb.Worksheets.Item("Landing Page")
$ws.Range("B8").Value = "What were the key themes from RBC's Q3 2025 earnings call?"
$ws.Range("B9").Value = "Summarize the management discussion from TD's Q2 2025 earnings call."

# --- "Popup Questions" sheet ---
$ws = $wb.Worksheets.Item("Popup Questions")
$ws.Range("C17").Value = "What outlook and guidance did RBC management provide in Q3 2025?"
$ws.Range("C18").Value = "What were the main topics analysts asked about in TD's Q2 2025 earnings call?"
$ws.Range("C19").Value = "What strategic initiatives did BMO management highlight in Q3 2025?"
$ws.Range("C20").Value = "Summarize the key points from Scotiabank's Q2 2025 management discussion."
$ws.Range("C21").Value = "Compare the key themes from RBC and TD's Q3 2025 earnings calls."

# --- "Dropdown - What is" sheet ---
$ws = $wb.Worksheets.Item("Dropdown - What is")
$ws.Range("B6").Value = "What is the key guidance RBC management provided in Q3 2025?"

# --- "Dropdown - Compare" sheet ---
$ws = $wb.Worksheets.Item("Dropdown - Compare")
$ws.Range("B3").Value = "Compare the management outlook from RBC and TD in Q3 2025."

# --- "Dropdown - How did" sheet ---
$ws = $wb.Worksheets.Item("Dropdown - How did")
$ws.Range("B2").Value = "How did TD management describe performance and outlook in Q2 2025?"
$ws.Range("B3").Value = "How did RBC management respond to analyst concerns in Q3 2025?"
$ws.Range("B4").Value = "How did BMO management describe their strategic priorities in Q3 2025?"

# --- "Dropdown - Summarize" sheet ---
$ws = $wb.Worksheets.Item("Dropdown - Summarize")
$ws.Range("B3").Value = "Summarize the key themes from Scotiabank's Q3 2025 earnings call."
$ws.Range("B6").Value = "Summarize the analyst Q&A session from RBC's Q3 2025 earnings call."
